$d = $word.ActiveDocument

# 1. Update the title text.
$d.Content.Find.Execute("Complex Test Document", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Test Document with Table", 2)

# 2. Update the intro paragraph text.
$d.Content.Find.Execute("This document has multiple tables.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This is a test document.", 2)

# 3. Remove the first (Key/Value/Status/Active) table entirely.
$d.Tables(1).Delete()

# 4. Remove the "Here is another table:" paragraph (text + its own mark).
$rng = $d.Content
$rng.Find.Execute("Here is another table:")
$rng.Paragraphs(1).Range.Delete()

# 5. Remove the last data row (Orange / 2.49 / 75) from the remaining table.
$d.Tables(1).Rows(4).Delete()

# 6. Apply the LightGrid-Accent1 table style.
$d.Tables(1).Style = "LightGrid-Accent1"

# 7. Update header row text.
$d.Content.Find.Execute("Product", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Name", 2)
$d.Content.Find.Execute("Price", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Age", 2)
$d.Content.Find.Execute("Stock", $true, $false, $false, $false, $false,
                         $true, 1, $false, "City", 2)

# 8. Update first data row text.
$d.Content.Find.Execute("Apple", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Alice", 2)
$d.Content.Find.Execute("1.99", $true, $false, $false, $false, $false,
                         $true, 1, $false, "30", 2)
$d.Content.Find.Execute("100", $true, $false, $false, $false, $false,
                         $true, 1, $false, "NYC", 2)

# 9. Update second data row text.
$d.Content.Find.Execute("Banana", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Bob", 2)
$d.Content.Find.Execute("0.99", $true, $false, $false, $false, $false,
                         $true, 1, $false, "25", 2)
$d.Content.Find.Execute("50", $true, $false, $false, $false, $false,
                         $true, 1, $false, "LA", 2)
